$p = $ppt.ActivePresentation

# 1. Update slide 1 title: "Review Activity Teams: Virtual" -> "Review Activity Teams"
$s1 = $p.Slides.Item(1)
$titleShp = $s1.Shapes.Item(1)
if ($titleShp.TextFrame.TextRange.Text -eq "Review Activity Teams: Virtual") {
    $titleShp.TextFrame.TextRange.Text = "Review Activity Teams"
}

# 2. Delete slide 2 ("Review Activity Teams: In-Person")
$s2 = $p.Slides.Item(2)
if ($s2.Shapes.Item(1).TextFrame.TextRange.Text -eq "Review Activity Teams: In-Person") {
    $s2.Delete()
}

# 3. Update the cached notes-slide slide-number field text "3" -> "2"
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    $np = $s.NotesPage
    for ($i = 1; $i -le $np.Shapes.Count; $i++) {
        $shp = $np.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "3") {
                $shp.TextFrame.TextRange.Text = "2"
            }
        }
    }
}

# 4. Update cached date field text across slide master / layouts / notes master
#    "10/24/2022" -> "2/3/2023"
#    "October 24, 2022" -> "February 3, 2023"
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $t = $shp.TextFrame.TextRange.Text
        if ($t -eq "10/24/2022") {
            $shp.TextFrame.TextRange.Text = "2/3/2023"
        } elseif ($t -eq "October 24, 2022") {
            $shp.TextFrame.TextRange.Text = "February 3, 2023"
        }
    }
}

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $lay = $layouts.Item($li)
    for ($i = 1; $i -le $lay.Shapes.Count; $i++) {
        $shp = $lay.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            $t = $shp.TextFrame.TextRange.Text
            if ($t -eq "10/24/2022") {
                $shp.TextFrame.TextRange.Text = "2/3/2023"
            } elseif ($t -eq "October 24, 2022") {
                $shp.TextFrame.TextRange.Text = "February 3, 2023"
            }
        }
    }
}

$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    $shp = $notesMaster.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $t = $shp.TextFrame.TextRange.Text
        if ($t -eq "10/24/2022") {
            $shp.TextFrame.TextRange.Text = "2/3/2023"
        } elseif ($t -eq "October 24, 2022") {
            $shp.TextFrame.TextRange.Text = "February 3, 2023"
        }
    }
}
